$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set column B width and style (target raw OOXML width is 15.375 chars;
# the COM ColumnWidth setter here quantizes to an MDW-7 pixel grid for this
# workbook's default font, so 14.7142855 is the closest input that lands on
# the nearest reachable stored width, 15.428571428571429)
$ws.Columns.Item(2).ColumnWidth = 14.7142855

# Apply text number format to column B (header + new cell) matching style s="1"
$ws.Range("B1:B2").NumberFormat = "@"

# Fill in new row 2 data (set C2 first so "qwe123" becomes shared string index 4,
# then B2 so "14412350000" becomes shared string index 5 - matches target order)
$ws.Range("C2").Value = "qwe123"
$ws.Range("B2").Value = "14412350000"

# Update selection
$ws.Range("C15").Select()

# Page setup
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
